# Update 苏州-漫展信息 workbook "想去人数" (F column) counts.
# These changes apply to both the "展览" sheet (rows as listed) and the
# "全部类型" sheet (which contains the same events shifted down by one row
# starting from row 24, due to an extra performance event inserted earlier).

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAll = $wb.Worksheets.Item("全部类型")

# Row -> (old, new) updates for "展览" sheet
$exhibitionUpdates = @{
    4  = 1272
    5  = 1120
    6  = 14155
    7  = 15826
    9  = 63
    18 = 83
    20 = 1227
    24 = 6255
    25 = 963
    26 = 1098
    27 = 5600
    30 = 129
    31 = 4586
}

foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Row -> new value updates for "全部类型" sheet (rows shifted by +1 from row 24 onward)
$allTypesUpdates = @{
    4  = 1272
    5  = 1120
    6  = 14155
    7  = 15826
    9  = 63
    18 = 83
    20 = 1227
    25 = 6255
    26 = 963
    27 = 1098
    28 = 5600
    31 = 129
    32 = 4586
}

foreach ($row in $allTypesUpdates.Keys) {
    $sheetAll.Range("F$row").Value = $allTypesUpdates[$row]
}
